$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.020.03"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "2.340.73"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'300.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").Value = "'99.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.62%  "
$ws.Range("D7").Value = "'0.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.77%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "'0.509"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.68%  "
$ws.Range("D10").Value = "'34.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.18%  "
$ws.Range("D11").Value = "'0.0787"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.38%  "
$ws.Range("D12").Value = "'7.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.87%  "
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "2.702.74"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "2.353.63"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "'13.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.80%  "
$ws.Range("D17").Value = "'0.805"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.92%  "
$ws.Range("D18").Value = "45.971.53"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").Value = "'12.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -9.16%  "
$ws.Range("D20").Value = "0.0₃0960"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").Value = "'5.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.55%  "
$ws.Range("D22").Value = "'66.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("D23").Value = "'243.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("D24").Value = "'2.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.40%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'1.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.18%  "
$ws.Range("D27").Value = "'40.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.72%  "
$ws.Range("D28").Value = "'2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.18%  "
$ws.Range("D29").Value = "'9.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.90%  "
$ws.Range("D30").Value = "'20.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").Value = "'3.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +14.58%  "
$ws.Range("E32").Value = "  +6.25%  "
$ws.Range("D33").Value = "'5.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.87%  "
$ws.Range("D34").Value = "'144.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").Value = "'0.0766"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.19%  "
$ws.Range("D36").Value = "'0.112"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("D37").Value = "'0.115"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.89%  "
$ws.Range("D38").Value = "'1.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("D39").Value = "'15.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.20%  "
$ws.Range("D40").Value = "'3.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.24%  "
$ws.Range("D41").Value = "'0.0297"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.68%  "
$ws.Range("D42").Value = "'3.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.69%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "1.845.51"
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("D45").Value = "'90.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("E46").Value = "  -7.74%  "
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").Value = "'70.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.43%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.184"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.22%  "
$ws.Range("D49").Value = "2.573.99"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'95.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.74%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'4.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
